# se actualiza el diccionario de datos
# Renombra el encabezado de la columna M (antes "NIT_Beneficiario") a
# "identificacion_Beneficiario" y ajusta el ancho de la columna a su nuevo
# contenido.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "identificacion_Beneficiario"
$ws.Columns.Item(13).ColumnWidth = 24.8

$ws.Range("M5").Select()
